# Applies the "adding averages and more checks" commit to the workbook.
#
# Summary of changes:
# 1. Header/title font: the bold header font loses its size-14 variant and the
#    bold "no-color" font gains a white font color (00FFFFFF) instead - i.e.
#    title (A1) and column header row (row 2) end up sharing one bold white font.
# 2. Training Dashboard sheet: H3:H9 ("PERIOD TO EXPIRE") values each drop by 8,
#    and I3:I9 ("LAST UPDATE") text changes from 08-Sep-2025 to 16-Sep-2025.
# 3. Exam Dashboard sheet: column E width grows from 10 to 15, and the COMMENTS
#    text in E3/E4 changes from "OK" to "date is valid".

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Training Dashboard
$ws2 = $wb.Worksheets.Item(2)   # Exam Dashboard

# --- 1. Font / style updates -------------------------------------------------
# Normalize the title font size down to the default (11pt) so it matches the
# header-row font, then paint both the title and the header rows white so
# they share a single bold white font (mirrors fonts count 3 -> 2 in styles.xml).
$ws1.Range("A1").Font.Size = 11
$ws1.Range("A1").Font.Color = 16777215
$ws1.Range("A2:K2").Font.Color = 16777215

$ws2.Range("A1").Font.Size = 11
$ws2.Range("A1").Font.Color = 16777215
$ws2.Range("A2:G2").Font.Color = 16777215

# --- 2. Training Dashboard data updates -------------------------------------
$ws1.Range("H3").Value = 334
$ws1.Range("H4").Value = 251
$ws1.Range("H5").Value = 647
$ws1.Range("H6").Value = 364
$ws1.Range("H7").Value = 423
$ws1.Range("H8").Value = 597
$ws1.Range("H9").Value = 181

# Keep these as literal text (not auto-converted into real Excel dates) by
# writing them with a leading text-qualifier apostrophe, same as typing a
# date-like value into a "General" formatted cell and forcing text mode.
$ws1.Range("I3").Value = "'16-Sep-2025"
$ws1.Range("I4").Value = "'16-Sep-2025"
$ws1.Range("I5").Value = "'16-Sep-2025"
$ws1.Range("I6").Value = "'16-Sep-2025"
$ws1.Range("I7").Value = "'16-Sep-2025"
$ws1.Range("I8").Value = "'16-Sep-2025"
$ws1.Range("I9").Value = "'16-Sep-2025"

# --- 3. Exam Dashboard updates -----------------------------------------------
# Column E width 10 -> 15 (ColumnWidth units run ~0.83 narrower than the raw
# OOXML column width, so 14.17 round-trips to an XML width of 15).
$ws2.Range("E1").EntireColumn.ColumnWidth = 14.17

$ws2.Range("E3").Value = "date is valid"
$ws2.Range("E4").Value = "date is valid"
